$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.312485933303833
$ws.Range("B1").Value = 3.609421253204346
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 3.083870649337769
$ws.Range("E1").Value = 2.627740621566772
